$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet entry (row 15): clocked in, not yet clocked out.
# The blank "Time out" (C15) makes the Delta/minutes/hours/money formulas
# on this row go negative (C15 - B15 with C15 = 0) -- this is the
# "new friction idea" from the commit message.
$ws.Range("A15").Value = 45578
$ws.Range("B15").Value = 0.42222222222222222

# Active selection moved to C15 (next cell the user would fill in)
$ws.Range("C15").Select()
